$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'60.239.60"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +3.27%  '
$ws.Range("D3").Value = "'3.213.80"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.03%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = "'539.09"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.35%  '
$ws.Range("D6").Value = "'146.65"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +4.92%  '
$ws.Range("D7").Value = "'1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("D8").Value = "'0.531"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +3.13%  '
$ws.Range("E10").Value = '  +4.06%  '
$ws.Range("D11").Value = "'0.435"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.08%  '
$ws.Range("D12").Value = "'3.770.46"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.12%  '
$ws.Range("E13").Value = '  -1.11%  '
$ws.Range("D14").Value = "'26.28"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.88%  '
$ws.Range("D15").Value = "'0.0000174"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.91%  '
$ws.Range("D16").Value = "'60.298.71"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +3.28%  '
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = "'3.209.79"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.74%  '
$ws.Range("B18").Value = 'Polkadot'
$ws.Range("C18").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D18").Value = "'6.28"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.78%  '
$ws.Range("D19").Value = "'13.20"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.49%  '
$ws.Range("E20").Value = '  +2.79%  '
$ws.Range("D21").Value = "'381.47"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.09%  '
$ws.Range("E22").Value = '  +0.08%  '
$ws.Range("D23").Value = "'0.528"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.18%  '
$ws.Range("D24").Value = "'70.22"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.52%  '
$ws.Range("D25").Value = "'8.91"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +11.23%  '
$ws.Range("E26").Value = '  +1.65%  '
$ws.Range("E27").Value = '  +0.16%  '
$ws.Range("D28").Value = "'0.0₃0909"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.28%  '
$ws.Range("D29").Value = "'6.23"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.86%  '
$ws.Range("E30").Value = '  +0.56%  '
$ws.Range("D31").Value = "'22.47"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.95%  '
$ws.Range("D32").Value = "'5.47"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +5.64%  '
$ws.Range("E33").Value = '  +3.72%  '
$ws.Range("D34").Value = "'6.65"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +6.32%  '
$ws.Range("D35").Value = "'157.18"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.63%  '
$ws.Range("D36").Value = "'1.36"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.73%  '
$ws.Range("B37").Value = 'Maker'
$ws.Range("C37").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D37").Value = "'2.798.76"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +5.56%  '
$ws.Range("B38").Value = 'EnergySwap'
$ws.Range("C38").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D38").Value = "'25.91"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.32%  '
$ws.Range("D39").Value = "'0.0709"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +4.18%  '
$ws.Range("E40").Value = '  +0.19%  '
$ws.Range("D41").Value = "'4.27"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.13%  '
$ws.Range("E42").Value = '  +3.94%  '
$ws.Range("D43").Value = "'0.722"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.67%  '
$ws.Range("D44").Value = "'0.0289"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +4.65%  '
$ws.Range("D45").Value = "'3.255.91"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.03%  '
$ws.Range("E46").Value = '  +2.88%  '
$ws.Range("D47").Value = "'0.103"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.56%  '
$ws.Range("D48").Value = "'6.19"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.63%  '
$ws.Range("D49").Value = "'0.806"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +7.14%  '
$ws.Range("D50").Value = "'20.82"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.83%  '
$ws.Range("D51").Value = "'273.45"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +10.63%  '
